$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values on row 3
$ws.Range("D3").Value = 0.001373819338487304
$ws.Range("E3").Value = 0.08924948771498312

# Update existing values on row 4
$ws.Range("B4").Value = 70
$ws.Range("C4").Value = 0.5544919454098156
$ws.Range("D4").Value = 0.006155956586378645
$ws.Range("E4").Value = 0.6059820549089501

# Update existing value on row 5
$ws.Range("D5").Value = 0.002743236579369405

# Add new row 6 with DWA data, copying the style of A5 for A6 (label column)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "DWA"

$ws.Range("B6").Value = 60
$ws.Range("C6").Value = 0.6
$ws.Range("D6").Value = 0.01285203008513388
$ws.Range("E6").Value = 0.5421821575419976
